# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Columns G (sexo), I (direccion-provincial-nombre) and J
# (edad-grupos-quinquenales) move from "dimension" to "measure":
#   - row 2: iaest-dimension:* / sdmx-dimension:* -> iaest-measure:*
#   - row 3: "dim" -> "medida"
#   - row 4: concept/URI datatype -> "xsd:int"
#   - row 5: their "mapping-*.xlsx" reference cells are removed entirely
#            (they no longer need an external code->concept mapping file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - iaest/sdmx identifier
$ws.Range("G2").Value = "iaest-measure:sexo"
$ws.Range("I2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("J2").Value = "iaest-measure:edad-grupos-quinquenales"

# Row 3 - "dim" becomes "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("J3").Value = "medida"

# Row 4 - datatype becomes "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"

# Row 5 - drop the now-unused mapping-file cells entirely
$ws.Range("G5").Clear()
$ws.Range("J5").Clear()
